$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New row 4: "Employment by industry" (APS source, same dates as row 2/3)
#    Insert above current row 4 ("Online job adverts...").
#    Template row-3's formatting (A=style1/left+center, B=style2), then
#    strip the style back off column A because the target row has no style
#    on A4.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial()
$ws.Cells.Item(4, 1).ClearFormats()

$ws.Cells.Item(4, 1).Value = "Employment by industry"
$ws.Cells.Item(4, 2).Value = "<a href='https://www.nomisweb.co.uk/datasets/apsnew'>Annual Population Survey</a>"
$ws.Cells.Item(4, 3).Value = "Jul 2021 - Jun 2022 (11/10/22)"
$ws.Cells.Item(4, 4).Value = "Oct 2021 - Sep 2022 (17/01/23)"

# ---------------------------------------------------------------------------
# After the above insert, the old rows 4-9 are now rows 5-10:
#   5 Online job adverts by local authority
#   6 Further education and skills achievements
#   7 Further education and skills achievements by sector subject area
#   8 Enterprise by employment size
#   9 Key Stage 4 (KS4) destinations
#  10 Key Stage 5 (KS5) destinations
#
# 2) New row 8: "Qualification by age and gender - NVQ" (APS source, new
#    date range). Insert above the current row 8 ("Enterprise by employment
#    size"), which carries no cell styling at all, and then paint column A
#    with the same format used for the other "A" header cells (style 1)
#    using a formats-only paste so no new style entries are created.
# ---------------------------------------------------------------------------
$ws.Rows.Item(8).Insert()
$ws.Range("A9:D9").Copy()
$ws.Range("A8:D8").PasteSpecial()
$ws.Range("A8:D8").ClearFormats()

$ws.Cells.Item(8, 1).Value = "Qualification by age and gender - NVQ"
$ws.Cells.Item(8, 2).Value = "<a href='https://www.nomisweb.co.uk/datasets/apsnew'>Annual Population Survey</a>"
$ws.Cells.Item(8, 3).Value = "Jan 2021 - Dec 2021 (20/04/21)"
$ws.Cells.Item(8, 4).Value = "Jan 2022 - Dec 2022 (19/04/22)"

$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(8, 1).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Row 9 is now "Enterprise by employment size" - refresh its release
#    dates from the Oct/Sept cycle to the new Mar cycle.
# ---------------------------------------------------------------------------
$ws.Cells.Item(9, 3).Value = "Mar 2022 (28/09/22)"
$ws.Cells.Item(9, 4).Value = "Mar 2023 (03/10/23)"

# ---------------------------------------------------------------------------
# 4) New row 10: "Enterprise by employment size and industry" - same source
#    and dates as row 9, inserted directly below it (above the KS4/KS5 rows
#    which currently sit at rows 10-11). Same no-style template + formats
#    paste trick as step 2.
# ---------------------------------------------------------------------------
$ws.Rows.Item(10).Insert()
$ws.Range("A11:D11").Copy()
$ws.Range("A10:D10").PasteSpecial()
$ws.Range("A10:D10").ClearFormats()

$ws.Cells.Item(10, 1).Value = "Enterprise by employment size and industry"
$ws.Cells.Item(10, 2).Value = "<a href='https://www.nomisweb.co.uk/datasets/idbrent'>ONS UK Business Count</a>"
$ws.Cells.Item(10, 3).Value = "Mar 2022 (28/09/22)"
$ws.Cells.Item(10, 4).Value = "Mar 2023 (03/10/23)"

$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(10, 1).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Column D needs to widen to fit the longer date strings now in it.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 26.55

# ---------------------------------------------------------------------------
# View state: zoom to 84%, and select A2:A12 (whole data column for column A).
# ---------------------------------------------------------------------------
$ws.Range("A2:A12").Select()
$excel.ActiveWindow.Zoom = 84
